$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing component position/rotation data (C14, C15) ---
$ws.Range("B15").Value = 121.98999999999999
$ws.Range("C15").Value = -113.2
$ws.Range("D15").Value = 0

$ws.Range("B16").Value = 115.79000000000001
$ws.Range("C16").Value = -115.95999999999999

# --- Insert two new rows for C20 / C21 right after C19 (before D1) ---
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = "C20"
$ws.Range("B21").Value = 123.34999999999999
$ws.Range("C21").Value = -80.260000000000005
$ws.Range("D21").Value = -90
$ws.Range("E21").Value = "top"

$ws.Range("A22").Value = "C21"
$ws.Range("B22").Value = 123.34999999999999
$ws.Range("C22").Value = -88.915000000000006
$ws.Range("D22").Value = 90
$ws.Range("E22").Value = "top"

# --- Insert one new row for D5 right after D4 (before J1) ---
# (row 27 currently holds J1, shifted down by the two rows inserted above)
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = "D5"
$ws.Range("B27").Value = 117.48
$ws.Range("C27").Value = -108.95999999999999
$ws.Range("D27").Value = -90
$ws.Range("E27").Value = "top"

# --- Update existing component position/rotation data further down the table ---
# (row numbers below already account for the three rows inserted above)
$ws.Range("B47").Value = 120.01000000000001
$ws.Range("C47").Value = -98.180000000000007

$ws.Range("B52").Value = 118.05
$ws.Range("C52").Value = -113.19
$ws.Range("D52").Value = 0

$ws.Range("B53").Value = 128.66
$ws.Range("C53").Value = -84.989999999999995

# --- Append new row for Y1 at the end of the table ---
$ws.Range("B53:D53").Copy($ws.Range("B54:D54"))
$ws.Range("A54").Value = "Y1"
$ws.Range("B54").Value = 123.34999999999999
$ws.Range("C54").Value = -84.159999999999997
$ws.Range("D54").Value = 90
$ws.Range("E54").Value = "top"
